$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-09-06 Saturday" "2025-09-07 Sunday"

Replace-Text "305÷9=" "760÷7="
Replace-Text "698÷2=" "683÷4="
Replace-Text "498÷3=" "163÷8="
Replace-Text "654÷4=" "664÷8="
Replace-Text "878÷6=" "196÷9="

Replace-Text "883÷4=" "980÷4="
Replace-Text "134÷9=" "164÷8="
Replace-Text "117÷3=" "873÷3="
Replace-Text "934÷4=" "746÷2="
Replace-Text "655÷5=" "390÷9="

Replace-Text "167÷3=" "570÷4="
Replace-Text "484÷8=" "959÷9="
Replace-Text "490÷4=" "927÷7="
Replace-Text "302÷4=" "358÷6="
Replace-Text "918÷3=" "828÷4="

Replace-Text "870÷9=" "741÷8="
Replace-Text "500÷2=" "432÷8="
Replace-Text "791÷7=" "534÷7="
Replace-Text "957÷4=" "820÷4="
Replace-Text "533÷4=" "542÷5="

Replace-Text "332÷6=" "910÷3="
Replace-Text "288÷2=" "486÷2="
Replace-Text "599÷7=" "872÷6="
Replace-Text "170÷8=" "710÷6="
Replace-Text "159÷3=" "545÷3="
